$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 126, which shifts the old row 126 down to row 127
$ws.Rows.Item(126).Insert()

# --- Update existing rows 116-125 with their new values ---

# Row 116
$ws.Range("D116").Value = 45223
$ws.Range("J116").Value = 50
$ws.Range("K116").Value = 10000
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = 10000
$ws.Range("P116").Value = 167

# Row 117
$ws.Range("D117").Value = 44754
$ws.Range("J117").Value = 60
$ws.Range("K117").Value = 11000
$ws.Range("L117").Value = 12000
$ws.Range("M117").Value = 11500
$ws.Range("P117").Value = 192

# Row 118
$ws.Range("I118").Value = "Primera"
$ws.Range("K118").Value = 8000
$ws.Range("L118").Value = 8000
$ws.Range("M118").Value = 8000
$ws.Range("N118").Value = "$/caja 60 unidades"
$ws.Range("P118").Value = 133
$ws.Range("Q118").Value = 60

# Row 119
$ws.Range("D119").Value = 44973
$ws.Range("I119").Value = "Segunda"
$ws.Range("J119").Value = 30
$ws.Range("K119").Value = 8500
$ws.Range("L119").Value = 8500
$ws.Range("M119").Value = 8500
$ws.Range("N119").Value = "$/caja 90 unidades"
$ws.Range("P119").Value = 94
$ws.Range("Q119").Value = 90

# Row 120
$ws.Range("D120").Value = 44999
$ws.Range("J120").Value = 50
$ws.Range("K120").Value = 10000
$ws.Range("L120").Value = 10000
$ws.Range("M120").Value = 10000
$ws.Range("O120").Value = "Región de Arica y Parinacota"
$ws.Range("P120").Value = 167

# Row 121
$ws.Range("D121").Value = 44608
$ws.Range("J121").Value = 100
$ws.Range("O121").Value = "Región Metropolitana"

# Row 122
$ws.Range("D122").Value = 44859
$ws.Range("J122").Value = 60
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 13000
$ws.Range("M122").Value = 12500
$ws.Range("P122").Value = 208

# Row 123
$ws.Range("D123").Value = 45063
$ws.Range("J123").Value = 50
$ws.Range("K123").Value = 8500
$ws.Range("L123").Value = 8500

# Row 124
$ws.Range("D124").Value = 45173
$ws.Range("J124").Value = 160
$ws.Range("K124").Value = 8000
$ws.Range("L124").Value = 9000
$ws.Range("M124").Value = 8500
$ws.Range("P124").Value = 142

# Row 125
$ws.Range("D125").Value = 44818
$ws.Range("J125").Value = 60
$ws.Range("K125").Value = 12000
$ws.Range("L125").Value = 13000
$ws.Range("M125").Value = 12500
$ws.Range("P125").Value = 208

# --- Row 126: brand new row created by the Insert() above, fill in all values ---
$ws.Range("A126").Value = 7
$ws.Range("B126").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C126").Value = "Ñuble"
$ws.Range("D126").Value = 45140
$ws.Range("E126").Value = 16
$ws.Range("F126").Value = 100112001
$ws.Range("G126").Value = "Berenjena"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 30
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = 9000
$ws.Range("N126").Value = "$/caja 60 unidades"
$ws.Range("O126").Value = "Región de Arica y Parinacota"
$ws.Range("P126").Value = 150
$ws.Range("Q126").Value = 60
$ws.Range("R126").Value = "Hortaliza"

# Apply the date number format (same as other D-column cells) to the new D126 cell
$ws.Range("D126").NumberFormat = $ws.Range("D125").NumberFormat
